$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A51").Value = "Davide Zeni "
$ws.Range("B51").Value = "Alessandro  Ruele | F.C. Gorillaz"
$ws.Range("C51").Value = "Federico Andreis | iMontagna"
$ws.Range("D51").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("E51").Value = "Eduardo  Grazioli  | FC Savignano"
$ws.Range("F51").Value = "Andrea Menolli | SdrumALA"
